# "Generate Report for Archive"
#
# The localization status report is refreshed: rows that were previously
# "Ready for handoff" have moved on to "In Translation". That text lives
# once in the shared-string table and is referenced from:
#   - Overview!E2, Overview!F2, Overview!E3, Overview!F3  (zh-cn / de-de status columns)
#   - zh-cn!C2,  zh-cn!C3   (Status column)
#   - de-de!C2,  de-de!C3   (Status column)
# Updating any one of the Range.Value assignments below updates the shared
# string used by every other cell that referenced it, exactly like Excel
# does when you retype the same cached string into a cell.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# The Status columns on the Overview sheet (E & F) were autosized to the
# "Ready for handoff" text; now that the shorter "In Translation" string is
# in place, re-autofit/narrow them to match (same columns on the zh-cn and
# de-de detail sheets hold the "Status" column, i.e. column C there).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
